$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 11:07"

# --- Row 5 (India) ---
$ws.Range("B5").Value = 5122846
$ws.Range("C5").Value = 6953
$ws.Range("E5").Value = 1014510
$ws.Range("G5").Value = 27
$ws.Range("H5").Value = 83257

# --- Row 24 (Filipinas) ---
$ws.Range("B24").Value = 276289
$ws.Range("C24").Value = 3375
$ws.Range("D24").Value = 208096
$ws.Range("E24").Value = 63408
$ws.Range("G24").Value = 53
$ws.Range("H24").Value = 4785

# --- Row 26 (Indonesia) ---
$ws.Range("B26").Value = 232628
$ws.Range("C26").Value = 3635
$ws.Range("D26").Value = 166686
$ws.Range("E26").Value = 56720
$ws.Range("G26").Value = 122
$ws.Range("H26").Value = 9222

# --- Rows 47/48: Polonia overtakes Japon ---
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 76571
$ws.Range("C47").Value = 837
$ws.Range("D47").Value = 62725
$ws.Range("E47").Value = 11593
$ws.Range("G47").Value = 16
$ws.Range("H47").Value = 2253

$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 76448
$ws.Range("D48").Value = 68532
$ws.Range("E48").Value = 6455
$ws.Range("H48").Value = 1461

# --- Row 57 (Singapur) ---
$ws.Range("B57").Value = 57532
$ws.Range("C57").Value = 18
$ws.Range("E57").Value = 550

# --- Row 67 (Afganistan) ---
$ws.Range("B67").Value = 38872
$ws.Range("C67").Value = 17
$ws.Range("D67").Value = 32505
$ws.Range("E67").Value = 4931

# --- Row 75 (El Salvador) ---
$ws.Range("B75").Value = 27249
$ws.Range("C75").Value = 86
$ws.Range("D75").Value = 20392
$ws.Range("E75").Value = 6056

# --- Rows 89/90: Croacia overtakes Grecia ---
$ws.Range("A89").Value = "Croacia"
$ws.Range("B89").Value = 14279
$ws.Range("C89").Value = 250
$ws.Range("D89").Value = 11933
$ws.Range("E89").Value = 2108
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 238

$ws.Range("A90").Value = "Grecia"
$ws.Range("B90").Value = 14041
$ws.Range("D90").Value = 3804
$ws.Range("E90").Value = 9921
$ws.Range("H90").Value = 316

# --- Row 102 (Finlandia) ---
$ws.Range("B102").Value = 8799
$ws.Range("C102").Value = 49
$ws.Range("E102").Value = 760

# --- Rows 110/111: Eslovaquia overtakes Mozambique ---
$ws.Range("A110").Value = "Eslovaquia"
$ws.Range("B110").Value = 6021
$ws.Range("C110").Value = 161
$ws.Range("D110").Value = 3288
$ws.Range("E110").Value = 2694
$ws.Range("G110").Value = 1

$ws.Range("A111").Value = "Mozambique"
$ws.Range("B111").Value = 5994
$ws.Range("D111").Value = 3267
$ws.Range("E111").Value = 2688
$ws.Range("H111").Value = 39

# --- Row 117 (Hong Kong) ---
$ws.Range("B117").Value = 4994
$ws.Range("C117").Value = 9
$ws.Range("D117").Value = 4682
$ws.Range("E117").Value = 210

# --- Row 126 (Eslovenia) ---
$ws.Range("B126").Value = 4058
$ws.Range("C126").Value = 104
$ws.Range("D126").Value = 2897
$ws.Range("E126").Value = 1025
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 136

# --- Rows 131/132: Lituania overtakes Tailandia ---
$ws.Range("A131").Value = "Lituania"
$ws.Range("B131").Value = 3504
$ws.Range("C131").Value = 62
$ws.Range("D131").Value = 2149
$ws.Range("E131").Value = 1268
$ws.Range("H131").Value = 87

$ws.Range("A132").Value = "Tailandia"
$ws.Range("B132").Value = 3490
$ws.Range("D132").Value = 3325
$ws.Range("E132").Value = 107
$ws.Range("H132").Value = 58

# --- Row 144 (Estonia) ---
$ws.Range("B144").Value = 2778
$ws.Range("C144").Value = 22
$ws.Range("D144").Value = 2337
$ws.Range("E144").Value = 377

# --- Row 193 (Brunei) ---
$ws.Range("D193").Value = 141
$ws.Range("E193").Value = 1
